$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 with latest transaction data
$ws.Range("A13").Value = 45708
$ws.Range("A13").NumberFormat = $ws.Range("A12").NumberFormat

$ws.Range("C13").Value = "Corte Adulto"
$ws.Range("D13").Value = 25
$ws.Range("E13").Value = "naty"
$ws.Range("G13").Value = "Efectivo"
